# fix error in excel formula
#
# accounting_example.xlsx had a few stale/incorrect formulas on Sheet1:
#   - I12:K12 carried a shared-formula template whose recorded range
#     included H12, even though H12 itself is a standalone formula. Re-enter
#     I12:K12 so the shared group is anchored correctly on I12.
#   - Row 20 ("Provisions") incorrectly added Net Income (row 19) on top of
#     Depreciation (row 17); it should just mirror Depreciation.
#   - Row 22 ("Free Cash Flow") needs to add back Net Income (row 19) after
#     netting Provisions (row 20) and Changes in Working Capital (row 21),
#     and the shared formula in row 22 needs to cover K22 too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mis-scoped shared formula on row 12 (Interest Revenue) ---
$ws.Range("I12:K12").Formula = "=I6*(0.1)*0.99"

# --- Row 20 (Provisions): should equal Depreciation only ---
$ws.Range("G20").Formula = "=G17"
$ws.Range("H20:K20").Formula = "=H17"

# --- Row 22 (Free Cash Flow): net out Provisions/Working-Capital change, add back Net Income ---
$ws.Range("G22").Formula = "=G20-G21+G19"
$ws.Range("H22:J22").Formula = "=H20-H21+H19"
$ws.Range("K22").Formula = "=K20-K21+K19"

# --- Restore the cursor/selection position recorded in the workbook ---
$ws.Range("L18").Select() | Out-Null
